$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Score"
$ws.Range("B1").Value = "Group"
$ws.Range("C1").Value = "Description"

$ws.Range("A1:C10").Select()
